$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet "Parametros" -> "Candidatas" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Candidatas"

# --- Add a new worksheet "Monitoreadas" right after the first one ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Monitoreadas"

# Header row re-using the same shared strings as sheet 1 (Bus k, Bus m, id, Rating)
$ws2.Range("B1").Value = "Bus k"
$ws2.Range("C1").Value = "Bus m"
$ws2.Range("D1").Value = "id"
$ws2.Range("E1").Value = "Rating"

# Copy the header formatting (bold font + border + centered alignment) from sheet1
$ws1.Range("B1").Copy()
$ws2.Range("B1:E1").PasteSpecial(-4122)

$data = @(
    @(111,113,1,2000),
    @(111,114,1,2000),
    @(112,113,1,2000),
    @(112,123,1,750),
    @(113,123,1,2000),
    @(113,215,1,2000),
    @(114,116,1,750),
    @(115,116,1,500),
    @(115,121,1,500),
    @(115,121,2,500),
    @(115,124,1,500),
    @(116,117,1,750),
    @(116,119,1,500),
    @(117,118,1,750),
    @(117,122,1,500),
    @(118,121,1,500),
    @(118,121,2,500),
    @(119,120,1,500),
    @(119,120,2,500),
    @(120,123,1,750),
    @(120,123,2,500),
    @(121,122,1,500),
    @(121,325,1,1000),
    @(123,217,1,750),
    @(211,213,1,500),
    @(211,214,1,500),
    @(212,213,1,500),
    @(212,223,1,500),
    @(213,223,1,500),
    @(214,216,1,500),
    @(215,216,1,750),
    @(215,221,1,450),
    @(215,221,2,450),
    @(215,224,1,500),
    @(216,217,1,750),
    @(216,219,1,750),
    @(217,218,1,500),
    @(217,222,1,500),
    @(218,221,1,500),
    @(218,221,2,500),
    @(219,220,1,750),
    @(219,220,2,500),
    @(220,223,1,750),
    @(220,223,2,500),
    @(221,222,1,500)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}

# Selection on the new sheet
$ws2.Range("B2:E46").Select()

# Re-select the original sheet's range last, so it remains the active/tabSelected sheet
$ws1.Range("I2:I11").Select()

$wb.Save()
